# Insert a new worksheet "pre_replication" between "cell_cycle" and "g1_to_s",
# listing the Reactome pre-replication complex proteins, and make it the
# active sheet (matching the commit: "create reactome pathway png for
# pre_replication proteins").

$wb = $excel.ActiveWorkbook

# Add the new sheet immediately before "g1_to_s" so the final order is
# cell_cycle, pre_replication, g1_to_s.
$g1ToS = $wb.Worksheets.Item("g1_to_s")
$newSheet = $wb.Worksheets.Add($g1ToS)
$newSheet.Name = "pre_replication"

# Column A: header + the pre-replication complex gene list.
$values = @(
    "pre_replication",
    "CDK2",
    "POLE2",
    "POLE",
    "POLE3",
    "POLE4",
    "DBF4",
    "CDC7",
    "ORC3",
    "ORC5",
    "ORC4",
    "ORC2",
    "MCM8",
    "ORC6",
    "ORC1",
    "CDC6",
    "MCM3",
    "MCM4",
    "MCM5",
    "MCM6",
    "MCM7",
    "MCM2",
    "MCM10",
    "CDC45",
    "RPA4",
    "RPA2",
    "RPA3",
    "RPA1",
    "GMNN",
    "CDT1",
    "POLA1",
    "POLA2",
    "PRIM1",
    "PRIM2"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $newSheet.Cells.Item($row, 1).Value = $values[$i]
}

# Matches the saved selection on the new sheet in the target workbook.
$newSheet.Range("B3").Select()
